$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 "62.620.22"
Set-TextValue 2 5 "  -2.06%  "
Set-TextValue 3 4 "3.200.24"
Set-TextValue 3 5 "  -3.16%  "
Set-TextValue 4 5 "  +0.02%  "
Set-TextValue 5 4 "595.03"
Set-TextValue 5 5 "  -1.02%  "
Set-TextValue 6 4 "135.97"
Set-TextValue 6 5 "  -5.19%  "
Set-TextValue 7 4 "1.00"
Set-TextValue 7 5 "  -0.01%  "
Set-TextValue 8 4 "3.197.66"
Set-TextValue 8 5 "  -3.17%  "
Set-TextValue 9 5 "  -2.83%  "
Set-TextValue 10 5 "  -3.46%  "
Set-TextValue 11 4 "5.34"
Set-TextValue 11 5 "  -2.25%  "
Set-TextValue 12 4 "0.455"
Set-TextValue 12 5 "  -3.36%  "
Set-TextValue 13 4 "0.0000239"
Set-TextValue 13 5 "  -4.22%  "
Set-TextValue 14 4 "33.61"
Set-TextValue 14 5 "  -3.55%  "
Set-TextValue 15 4 "3.724.31"
Set-TextValue 15 5 "  -3.17%  "
Set-TextValue 16 5 "  -0.29%  "
Set-TextValue 17 4 "3.201.20"
Set-TextValue 17 5 "  -3.14%  "
Set-TextValue 18 4 "62.699.62"
Set-TextValue 18 5 "  -2.05%  "
Set-TextValue 19 4 "6.71"
Set-TextValue 19 5 "  -2.68%  "
Set-TextValue 20 4 "462.85"
Set-TextValue 20 5 "  -3.85%  "
Set-TextValue 21 4 "14.03"
Set-TextValue 21 5 "  -2.06%  "
Set-TextValue 22 4 "0.713"
Set-TextValue 22 5 "  -4.04%  "
Set-TextValue 23 4 "7.69"
Set-TextValue 23 5 "  -4.09%  "
Set-TextValue 24 4 "13.57"
Set-TextValue 24 5 "  +0.46%  "
Set-TextValue 25 4 "83.84"
Set-TextValue 25 5 "  -0.52%  "
Set-TextValue 26 5 "  -0.11%  "
Set-TextValue 27 4 "2.72"
Set-TextValue 27 5 "  -2.23%  "
Set-TextValue 28 5 "  -0.08%  "
Set-TextValue 29 4 "7.94"
Set-TextValue 29 5 "  -3.85%  "
Set-TextValue 30 4 "6.93"
Set-TextValue 30 5 "  -5.72%  "
Set-TextValue 31 4 "2.08"
Set-TextValue 31 5 "  -3.33%  "
Set-TextValue 32 4 "27.52"
Set-TextValue 32 5 "  -3.26%  "
Set-TextValue 33 4 "0.102"
Set-TextValue 33 5 "  -4.00%  "
Set-TextValue 34 5 "  -4.03%  "
Set-TextValue 35 5 "  -4.67%  "
Set-TextValue 36 4 "5.88"
Set-TextValue 36 5 "  -1.89%  "
Set-TextValue 37 4 "51.60"
Set-TextValue 37 5 "  -3.24%  "
Set-TextValue 38 4 "0.0₃0697"
Set-TextValue 38 5 "  -7.57%  "
Set-TextValue 39 5 "  -1.09%  "
Set-TextValue 40 2 "Maker"
Set-TextValue 40 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 40 4 "3.013.49"
Set-TextValue 40 5 "  -0.30%  "
Set-TextValue 41 2 "Bittensor"
Set-TextValue 41 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 41 4 "417.04"
Set-TextValue 41 5 "  -2.90%  "
Set-TextValue 42 5 "  +5.03%  "
Set-TextValue 43 4 "8.12"
Set-TextValue 43 5 "  -3.74%  "
Set-TextValue 44 4 "2.64"
Set-TextValue 44 5 "  -5.14%  "
Set-TextValue 45 4 "0.254"
Set-TextValue 45 5 "  -5.63%  "
Set-TextValue 46 5 "  -2.43%  "
Set-TextValue 47 4 "36.05"
Set-TextValue 47 5 "  +0.45%  "
Set-TextValue 48 2 "USDe"
Set-TextValue 48 3 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue 48 4 "0.999"
Set-TextValue 48 5 "  -0.08%  "
Set-TextValue 49 2 "InjectiveProtocol"
Set-TextValue 49 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 49 4 "26.01"
Set-TextValue 49 5 "  -1.40%  "
Set-TextValue 50 4 "2.31"
Set-TextValue 50 5 "  -0.63%  "
Set-TextValue 51 4 "125.20"
Set-TextValue 51 5 "  +0.66%  "
